$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title placeholder (shape 1 / id=2, "Titre 1") -------------------------
$title = $s.Shapes.Item(1)

# Reposition / resize (matches target <a:xfrm> off/ext in EMU)
$title.Left = 120
$title.Top = 69.42764
$title.Width = 720
$title.Height = 135.20387

$tf = $title.TextFrame
$tf.AutoSize = 0
$tr = $tf.TextRange
$tr.Text = "Présentation revue 0"
$tr.Font.Size = 72
$tr.Font.Bold = $true
$tr.Font.Italic = $true
$tr.Font.Underline = $true

# --- Subtitle placeholder (shape 2 / id=3, "Sous-titre 2") ------------------
# Removed entirely in the target deck.
$s.Shapes.Item(2).Delete()
